# Apply updated betting-odds values to row 5 and row 8 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("G5").Value  = 4.9
$ws.Range("H5").Value  = 3.4
$ws.Range("J5").Value  = 5.2
$ws.Range("K5").Value  = 2.12
$ws.Range("L5").Value  = 2.2
$ws.Range("Q5").Value  = 2.02
$ws.Range("W5").Value  = 11
$ws.Range("X5").Value  = 27
$ws.Range("Y5").Value  = 17
$ws.Range("AA5").Value = 60
$ws.Range("AB5").Value = 70
$ws.Range("AC5").Value = 8
$ws.Range("AD5").Value = 6.8
$ws.Range("AE5").Value = 18.5
$ws.Range("AF5").Value = 110
$ws.Range("AG5").Value = 5.8
$ws.Range("AH5").Value = 7.2
$ws.Range("AK5").Value = 14.5
$ws.Range("AO5").Value = 30
$ws.Range("AP5").Value = 37
$ws.Range("AQ5").Value = 200
$ws.Range("AR5").Value = 250
$ws.Range("AS5").Value = 500
$ws.Range("AT5").Value = 2.5
$ws.Range("AU5").Value = 7.7
$ws.Range("AX5").Value = 7.9
$ws.Range("AY5").Value = 18
$ws.Range("AZ5").Value = 26
$ws.Range("BA5").Value = 60

# Row 8 updates
$ws.Range("Q8").Value = 2.35
$ws.Range("R8").Value = 1.57
